{"js": "  // Replace each old two-digit-division answer string with its new value.\n  // Mapping reflects the commit diff, applied in document order so the\n  // one coincidental value reuse (row 12 old text == row 25 new text)\n  // never double-matches.\n  const replacements = [\n    [\"49\u00f72=24, 1\", \"19\u00f74=4, 3\"],\n    [\"80\u00f79=8, 8\", \"42\u00f72=21, 0\"],\n    [\"21\u00f72=10, 1\", \"90\u00f78=11, 2\"],\n    [\"56\u00f75=11, 1\", \"32\u00f79=3, 5\"],\n    [\"39\u00f78=4, 7\", \"63\u00f74=15, 3\"],\n    [\"92\u00f75=18, 2\", \"47\u00f75=9, 2\"],\n    [\"66\u00f73=22, 0\", \"58\u00f72=29, 0\"],\n    [\"87\u00f75=17, 2\", \"28\u00f72=14, 0\"],\n    [\"69\u00f75=13, 4\", \"61\u00f75=12, 1\"],\n    [\"38\u00f73=12, 2\", \"34\u00f78=4, 2\"],\n    [\"50\u00f73=16, 2\", \"89\u00f72=44, 1\"],\n    [\"28\u00f75=5, 3\", \"72\u00f74=18, 0\"],\n    [\"78\u00f79=8, 6\", \"99\u00f74=24, 3\"],\n    [\"83\u00f77=11, 6\", \"53\u00f78=6, 5\"],\n    [\"28\u00f79=3, 1\", \"64\u00f73=21, 1\"],\n    [\"24\u00f76=4, 0\", \"42\u00f78=5, 2\"],\n    [\"45\u00f73=15, 0\", \"27\u00f76=4, 3\"],\n    [\"70\u00f72=35, 0\", \"85\u00f78=10, 5\"],\n    [\"14\u00f78=1, 6\", \"67\u00f75=13, 2\"],\n    [\"33\u00f75=6, 3\", \"86\u00f74=21, 2\"],\n    [\"77\u00f77=11, 0\", \"21\u00f73=7, 0\"],\n    [\"94\u00f77=13, 3\", \"32\u00f75=6, 2\"],\n    [\"72\u00f75=14, 2\", \"44\u00f76=7, 2\"],\n    [\"12\u00f74=3, 0\", \"62\u00f76=10, 2\"],\n    [\"49\u00f74=12, 1\", \"28\u00f75=5, 3\"],\n  ];\n\n  const body = context.document.body;\n  for (const [oldText, newText] of replacements) {\n    const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    found.load(\"items\");\n    await context.sync();\n\n    if (found.items.length === 0) {\n      throw new Error(`Text not found: ${oldText}`);\n    }\n\n    for (const range of found.items) {\n      range.insertText(newText, \"Replace\");\n    }\n    await context.sync();\n  }\n", "ps1": "# Update each two-digit-division answer cell to its new value, per the\n# commit diff. Uses Range.Find/Execute (wdReplaceAll = 2) scoped to the\n# whole document body. Pairs are applied in document order so the one\n# coincidental value reuse (row 12's old text equals row 25's new text)\n# never causes a find to match a freshly-written replacement.\n\n$replacements = @(\n    @(\"49\u00f72=24, 1\", \"19\u00f74=4, 3\"),\n    @(\"80\u00f79=8, 8\", \"42\u00f72=21, 0\"),\n    @(\"21\u00f72=10, 1\", \"90\u00f78=11, 2\"),\n    @(\"56\u00f75=11, 1\", \"32\u00f79=3, 5\"),\n    @(\"39\u00f78=4, 7\", \"63\u00f74=15, 3\"),\n    @(\"92\u00f75=18, 2\", \"47\u00f75=9, 2\"),\n    @(\"66\u00f73=22, 0\", \"58\u00f72=29, 0\"),\n    @(\"87\u00f75=17, 2\", \"28\u00f72=14, 0\"),\n    @(\"69\u00f75=13, 4\", \"61\u00f75=12, 1\"),\n    @(\"38\u00f73=12, 2\", \"34\u00f78=4, 2\"),\n    @(\"50\u00f73=16, 2\", \"89\u00f72=44, 1\"),\n    @(\"28\u00f75=5, 3\", \"72\u00f74=18, 0\"),\n    @(\"78\u00f79=8, 6\", \"99\u00f74=24, 3\"),\n    @(\"83\u00f77=11, 6\", \"53\u00f78=6, 5\"),\n    @(\"28\u00f79=3, 1\", \"64\u00f73=21, 1\"),\n    @(\"24\u00f76=4, 0\", \"42\u00f78=5, 2\"),\n    @(\"45\u00f73=15, 0\", \"27\u00f76=4, 3\"),\n    @(\"70\u00f72=35, 0\", \"85\u00f78=10, 5\"),\n    @(\"14\u00f78=1, 6\", \"67\u00f75=13, 2\"),\n    @(\"33\u00f75=6, 3\", \"86\u00f74=21, 2\"),\n    @(\"77\u00f77=11, 0\", \"21\u00f73=7, 0\"),\n    @(\"94\u00f77=13, 3\", \"32\u00f75=6, 2\"),\n    @(\"72\u00f75=14, 2\", \"44\u00f76=7, 2\"),\n    @(\"12\u00f74=3, 0\", \"62\u00f76=10, 2\"),\n    @(\"49\u00f74=12, 1\", \"28\u00f75=5, 3\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\nWrite-Output \"Replaced $($replacements.Count) answers.\"\n"}
